# Update cryptocurrency price/volume table to reflect the latest scrape.
# Generated from the OOXML diff: updates D (Price) / E (Volume 1h) values for
# existing rows, and replaces rows 48-50 (BabyDogeCoin dropped out of the top
# list; Algorand/Cronos shift up and EnergySwap is newly listed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.622.04'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.596.67'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.43'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0614'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.54'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = '1.822.71'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '1.597.80'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.65'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = '26.616.21'
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.53'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("E23").Value = '  -3.21%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.68'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.13'
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.26'
$ws.Range("E34").Value = '  +17.95%  '
$ws.Range("D35").Value = '1.277.96'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.49'
$ws.Range("E36").Value = '  +0.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.597'
$ws.Range("E38").Value = '  -3.79%  '
$ws.Range("E39").Value = '  -2.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.822'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.771'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.61'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '1.734.20'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.57'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.102'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0513'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.48'
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("E51").Value = '  +0.11%  '
